$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.741.43"
$ws.Range("E2").Value = "  +2.58%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.718.64"
$ws.Range("E3").Value = "  +2.49%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.994"
$ws.Range("E4").Value = "  -0.69%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.35"
$ws.Range("E5").Value = "  +0.94%  "

# Row 6
$ws.Range("E6").Value = "  +0.73%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").Value = "  -0.78%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.01"
$ws.Range("E8").Value = "  +12.05%  "

# Row 9
$ws.Range("E9").Value = "  +4.33%  "

# Row 10
$ws.Range("E10").Value = "  +0.95%  "

# Row 11
$ws.Range("E11").Value = "  +0.65%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.960.43"
$ws.Range("E12").Value = "  +2.44%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.706.57"
$ws.Range("E13").Value = "  +1.59%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.24"
$ws.Range("E14").Value = "  +2.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.566"
$ws.Range("E15").Value = "  +5.86%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.82"
$ws.Range("E16").Value = "  +2.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.728.51"
$ws.Range("E17").Value = "  +2.54%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.92"
$ws.Range("E18").Value = "  +2.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.96"
$ws.Range("E19").Value = "  -2.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0748"
$ws.Range("E20").Value = "  +1.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.992"
$ws.Range("E21").Value = "  -0.84%  "

# Row 22
$ws.Range("E22").Value = "  +3.37%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.69"
$ws.Range("E23").Value = "  +4.88%  "

# Row 24
$ws.Range("E24").Value = "  +0.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.42"
$ws.Range("E25").Value = "  +0.68%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.56"
$ws.Range("E26").Value = "  +3.91%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.60"
$ws.Range("E27").Value = "  +0.75%  "

# Row 28
$ws.Range("E28").Value = "  +0.99%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.994"
$ws.Range("E29").Value = "  -0.74%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  +1.38%  "

# Row 31
$ws.Range("E31").Value = "  +1.03%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.43"
$ws.Range("E32").Value = "  +1.64%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.556.95"
$ws.Range("E33").Value = "  +1.00%  "

# Row 34
$ws.Range("E34").Value = "  +4.13%  "

# Row 35
$ws.Range("E35").Value = "  -1.95%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.966"
$ws.Range("E36").Value = "  +5.53%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.614"
$ws.Range("E37").Value = "  +4.25%  "

# Row 38
$ws.Range("E38").Value = "  +0.38%  "

# Row 39
$ws.Range("E39").Value = "  +0.03%  "

# Row 40
$ws.Range("E40").Value = "  +1.78%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.21"
$ws.Range("E41").Value = "  +4.91%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.87"
$ws.Range("E42").Value = "  +5.94%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.993"
$ws.Range("E43").Value = "  -0.80%  "

# Row 44
$ws.Range("E44").Value = "  +1.19%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.862.20"
$ws.Range("E45").Value = "  +2.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.791"
$ws.Range("E46").Value = "  +1.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.69"
$ws.Range("E47").Value = "  +9.80%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "91.48"
$ws.Range("E48").Value = "  +1.24%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.106"
$ws.Range("E49").Value = "  +2.65%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.29"
$ws.Range("E50").Value = "  +3.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0104"
$ws.Range("E51").Value = "  -3.37%  "
